$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = 'na_rate_threshold-0.1'

$col1 = @(
    'AMS3000',
    'AMS3001',
    'AMS3002',
    'AMS3008',
    'AMS3009',
    'AMS3010',
    'AMS3011',
    'AMS3012',
    'AMS3013',
    'AMS3014',
    'AMS3015',
    'AMS3016',
    'AMS3017',
    'AMS3018',
    'AMS3019',
    'AMS3020',
    'AMS3021',
    'AMS3022',
    'AMS3023',
    'AMS3024',
    'AMS3025',
    'AMS3026',
    'AMS3027',
    'AMS3028',
    'AMS3029',
    'AMS3030',
    'AMS3031',
    'AMS3117',
    'AMS3119',
    'AMS3123',
    'AMS3126',
    'AMS3131',
    'AMS3132',
    'AMS3138',
    'AMS3140',
    'AMS3141',
    'AMS3143',
    'AMS3144',
    'AMS3145',
    'AMS3147',
    'AMS3150',
    'AMS3152',
    'AMS3153',
    'AMS3158',
    'AMS3160',
    'AMS3161',
    'AMS3162',
    'AMS3163',
    'AMS3164',
    'AMS3165',
    'AMS3166',
    'AMS3167',
    'AMS3169',
    'AMS3174',
    'AMS3181',
    'AMS3183',
    'AMS3184',
    'AMS3185',
    'AMS3188',
    'AMS3189',
    'AMS3190',
    'AMS3191',
    'AMS3195',
    'AMS3202',
    'AMS3203',
    'AMS3204',
    'AMS3205',
    'AMS3206',
    'AMS3207',
    'AMS3219',
    'AMS3220',
    'AMS3223',
    'AMS3228',
    'AMS3231',
    'AMS3234',
    'AMS3245',
    'AMS3272',
    'AMS3276',
    'AMS3285',
    'AMS3292',
    'AMS3296',
    'AMS3311',
    'AMS3315',
    'AMS3322',
    'AMS3326',
    'AMS3332',
    'AMS3335',
    'AMS3348',
    'AMS3351',
    'AMS3353',
    'AMS3357',
    'AMS3359',
    'AMS3363',
    'AMS3365',
    'AMS3371',
    'AMS3373',
    'AMS3377',
    'AMS3382',
    'AMS3385',
    'AMS3387',
    'AMS3391',
    'AMS3394',
    'AMS3396',
    'AMS3400',
    'AMS3403',
    'AMS3405',
    'AMS3409',
    'AMS3412',
    'AMS3414',
    'AMS3419',
    'AMS3422',
    'AMS3423',
    'AMS3425',
    'AMS3430',
    'AMS3433',
    'AMS3434',
    'AMS3436',
    'AMS3441',
    'AMS3444',
    'AMS3445',
    'AMS3447',
    'AMS3452',
    'AMS3455',
    'AMS3456',
    'AMS3458',
    'AMS3536',
    'AMS3539',
    'AMS3545',
    'AMS3547',
    'AMS3550',
    'AMS3556',
    'AMS3558',
    'AMS3561',
    'AMS3567',
    'AMS3569',
    'AMS3572',
    'AMS3578',
    'AMS3579',
    'AMS3580',
    'AMS3585',
    'AMS3589',
    'AMS3593',
    'AMS3596',
    'AMS3598',
    'AMS3600',
    'AMS3602',
    'AMS3607',
    'AMS3609',
    'AMS3610',
    'AMS3611',
    'AMS3618',
    'AMS3622',
    'AMS3624',
    'AMS3629',
    'AMS3633',
    'AMS3723',
    'AMS3725',
    'AMS3726',
    'AMS3727',
    'AMS3728',
    'AMS3729',
    'AMS3744',
    'AMS3745',
    'AMS3749',
    'AMS3753',
    'AMS3756',
    'AMS3759',
    'AMS3760',
    'AMS3761',
    'AMS3762',
    'AMS3797',
    'AMS3798',
    'AMS3799',
    'AMS3800',
    'AMS3801',
    'AMS3802',
    'AMS3803',
    'AMS3804',
    'AMS3805',
    'AMS3807',
    'AMS3808',
    'AMS3809',
    'AMS3812',
    'AMS3813',
    'AMS3814',
    'AMS3815',
    'AMS3830',
    'AMS3834',
    'AMS3847',
    'AMS3848',
    'AMS3849',
    'AMS3850',
    'AMS3851',
    'AMS3853',
    'AMS3854',
    'AMS3855',
    'AMS3857',
    'AMS3858',
    'AMS3859',
    'AMS3860',
    'AMS3861',
    'AMS3866',
    'AMS3867',
    'AMS3868',
    'AMS3876',
    'AMS3877',
    'AMS3878',
    'AMS3882',
    'AMS3883',
    'AMS3884',
    'AMS3890',
    'AMS3892',
    'AMS3897',
    'AMS3898',
    'AMS3899',
    'AMS3911',
    'AMS3912',
    'AMS3913',
    'AMS3923',
    'AMS3924',
    'AMS3925',
    'AMS3929',
    'AMS3930',
    'AMS3931',
    'AMS3939',
    'AMS3944',
    'AMS3945',
    'AMS3946',
    'AMS3956',
    'AMS3957',
    'AMS3958',
    'AMS3985',
    'AMS3990',
    'AMS3993',
    'AMS3994',
    'AMS3995'
)

$col2 = @(
    '# Inq w/in 3 M',
    '# Inq w/in 12 M',
    '# Inq w/in 24 M',
    '# Non-Util Inq w/in 1 Month',
    '# Non-Util Inq w/in 3 M',
    '# Non-Util Inq w/in 12 M',
    '# Non-Util Inq w/in 24 M',
    '% Inq w/in 1 Month to Inq w/in 12 M',
    '% Inq w/in 1 Month to Inq w/in 24 M',
    '% Util Inq w/in 3 M to Inq w/in 12 M',
    '% Util Inq w/in 3 M to Inq w/in 24 M',
    '% Util Inq w/in 12 M to Inq w/in 24 M',
    '% Util Inq w/in 1 Month to Inq w/in 12 M',
    '% Util Inq w/in 1 Month to Inq w/in 24 M',
    '% Non-Util Inq w/in 3 M to Inq w/in 12 M',
    '% Non-Util Inq w/in 3 M to Inq w/in 24 M',
    '% Non-Util Inq w/in 12 M to Inq w/in 24 M',
    '% Non-Util Inq w/in 1 Month to Inq w/in 12 M',
    '% Non-Util Inq w/in 1 Month to Inq w/in 24 M',
    '# Inq w/in 6 M',
    '# Util Inq w/in 6 M',
    '# Non-Util Inq w/in 6 M',
    '# Inq w/in 1 Month',
    '# Util Inq w/in 1 Month',
    '# Util Inq w/in 3 M',
    '# Util Inq w/in 12 M',
    '# Util Inq w/in 24 M',
    'Age Oldest Mortgage Trade',
    'Age Oldest Retail Trade',
    'Age Newest Auto Trade',
    'Age Newest Dept Store Trades',
    'Age Newest Retail Trade',
    'Age Newest Sales Finance Trade',
    '# Open Auto Trades',
    '# Open Credit Union Trades',
    '# Open Dept Store Trades',
    '# Open Mortgage Trades',
    '# Open PerFin and StuLoan Trades',
    '# Open Retail Trades',
    '# Open Sales Finance Trades',
    '# Bkcrd Trades w/Update w/in 3 M w/ Bal > $0 ',
    '# Dept Store Trades w/Update w/in 3 M w/ Bal > $0 ',
    '# Inst Trades w/Update w/in 3 M w/ Bal > $0 ',
    '# Retail Trades w/Update w/in 3 M w/ Bal > $0 ',
    'T Bal Open Auto Trades w/Update w/in 3 M',
    'T Bal Open Bkcrd Trades w/Update w/in 3 M',
    'T Bal Open Credit Union Trades w/Update w/in 3 M',
    'T Bal Open Dept Store Trades w/Update w/in 3 M',
    'T Bal Open Inst Trades w/Update w/in 3 M',
    'T Bal Open Mortgage Trades w/Update w/in 3 M',
    'T Bal Open PerFin and StuLoan Trades w/Update w/in 3 M',
    'T Bal Open Retail Trades w/Update w/in 3 M',
    'T Bal Open Sales Finance Trades w/Update w/in 3 M',
    'T Bal Open StuLoan Trades w/Update w/in 3 M',
    'T Bal Closed Trades w/Update w/in 3 M',
    'T Bal Closed Bkcrd Trades w/Update w/in 3 M',
    'T Bal Closed Credit Union Trades w/Update w/in 3 M',
    'T Bal Closed Dept Store Trades w/Update w/in 3 M',
    'T Bal Closed PerFin and StuLoan Trades w/Update w/in 3 M',
    'T Bal Closed Retail Trades w/Update w/in 3 M',
    'T Bal Closed Rev Trades w/Update w/in 3 M',
    'T Bal Closed Sales Finance Trades w/Update w/in 3 M',
    'Age Newest Mortgage Trade',
    'T Loan Amount Open StuLoan Trades w/Update w/in 3 M',
    'T Loan Amount Open Auto Trades w/Update w/in 3 M',
    'T High Credit Open Bkcrd Trades w/Update w/in 3 M',
    'T High Credit Open Dept Store Trades w/Update w/in 3 M',
    'T Loan Amount Open Inst Trades w/Update w/in 3 M',
    'T Loan Amount Open Mortgage Trades w/Update w/in 3 M',
    '# Dept Store Trades w/ PD > $0',
    'T High Credit Open Retail Trades w/Update w/in 3 M',
    '# Retail Trades w/ PD > $0',
    'T PD  Bkcrd Trades w/Update w/in 3 M',
    'T PD Inst Trades w/Update w/in 3 M',
    'T PD Retail Trades w/Update w/in 3 M',
    'T PD Retail Trades ',
    '# 30 DPD Occur w/in 24 M Dept Store Trades',
    '# 30 DPD Occur w/in 24 M Retail Trades',
    '# 60 DPD Occur w/in 12 M Retail Trades',
    '# 60 DPD Occur w/in 24 M Dept Store Trades',
    '# 60 DPD Occur w/in 24 M Retail Trades',
    '# 90 DPD Occur w/in 24 M Dept Store Trades',
    '# 90 DPD Occur w/in 24 M Retail Trades',
    '# 120-180 or More DPD Occur w/in 24 M Dept Store Trades',
    '# 120-180 or More DPD Occur w/in 24 M Retail Trades',
    '# Dept Store Trades Satis w/in 3 M',
    '# Retail Trades Satis w/in 3 M',
    '# Retail Trades Satis w/in 6 M',
    '# 60-180 or More DPD Occur w/in 24 M Dept Store Trades',
    '# 60-180 or More DPD Occur w/in 24 M Retail Trades',
    '# 90-180 or More DPD Occur w/in 24 M Dept Store Trades',
    '# 90-180 or More DPD Occur w/in 24 M Retail Trades',
    '# 60-180 or More DPD Occur w/in 12 M Dept Store Trades',
    '# 60-180 or More DPD Occur w/in 12 M Retail Trades',
    '# 90-180 or More DPD Occur w/in 12 M Dept Store Trades',
    '# 90-180 or More DPD Occur w/in 12 M Retail Trades',
    '# Retail Trades Always Satis',
    '# Dept Store Tr Wo Rat 30 DPD w/in 3 M',
    '# Retail Tr Wo Rat 30 DPD w/in 3 M',
    '# Sales Finance Tr Wo Rat 30 DPD w/in 3 M',
    '# Dept Store Tr Wo Rat 60 DPD w/in 3 M',
    '# Retail Tr Wo Rat 60 DPD w/in 3 M',
    '# Sales Finance Tr Wo Rat 60 DPD w/in 3 M',
    '# Dept Store Tr Wo Rat 90 DPD w/in 3 M',
    '# Retail Tr Wo Rat 90 DPD w/in 3 M',
    '# Sales Finance Tr Wo Rat 90 DPD w/in 3 M',
    '# Dept Store Tr Wo Rat 120-180 or More DPD w/in 3 M',
    '# Retail Tr Wo Rat 120-180 or More DPD w/in 3 M',
    '# Sales Finance Tr Wo Rat 120-180 or More DPD w/in 3 M',
    '# Dept Store Tr Wo Rat 30 DPD w/in 6 M',
    '# PerFin and StuLoan Tr Wo Rat 30 DPD w/in 6 M',
    '# Retail Tr Wo Rat 30 DPD w/in 6 M',
    '# Sales Finance Tr Wo Rat 30 DPD w/in 6 M',
    '# Dept Store Tr Wo Rat 60 DPD w/in 6 M',
    '# PerFin and StuLoan Tr Wo Rat 60 DPD w/in 6 M',
    '# Retail Tr Wo Rat 60 DPD w/in 6 M',
    '# Sales Finance Tr Wo Rat 60 DPD w/in 6 M',
    '# Dept Store Tr Wo Rat 90 DPD w/in 6 M',
    '# PerFin and StuLoan Tr Wo Rat 90 DPD w/in 6 M',
    '# Retail Tr Wo Rat 90 DPD w/in 6 M',
    '# Sales Finance Tr Wo Rat 90 DPD w/in 6 M',
    '# Dept Store Tr Wo Rat 120-180 or More DPD w/in 6 M',
    '# PerFin and StuLoan Tr Wo Rat 120-180 or More DPD w/in 6 M',
    '# Retail Tr Wo Rat 120-180 or More DPD w/in 6 M',
    '# Sales Finance Tr Wo Rat 120-180 or More DPD w/in 6 M',
    '# Auto Tr Wo Rat Ever 30 DPD',
    '# Dept Store Tr Wo Rat Ever 30 DPD',
    '# Retail Tr Wo Rat Ever 30 DPD',
    '# Auto Tr Wo Rat Ever 60 DPD',
    '# Dept Store Tr Wo Rat Ever 60 DPD',
    '# Retail Tr Wo Rat Ever 60 DPD',
    '# Auto Tr Wo Rat Ever 90 DPD',
    '# Dept Store Tr Wo Rat Ever 90 DPD',
    '# Retail Tr Wo Rat Ever 90 DPD',
    '# Auto Tr Wo Rat Ever 120-180 or More DPD',
    '# Dept Store Tr Wo Rat Ever 120-180 or More DPD',
    '# Retail Tr Wo Rat Ever120-180 or More DPD',
    'Worst Rating w/in 3 M Auto Trades',
    'Worst Rating w/in 3 M Mortgage Trades',
    '# Dept Store Trades w/ MajDerog Reported w/in 6 M',
    '# Retail Trades w/ MajDerog Reported w/in 6 M',
    '# Auto Trades w/ MajDerog Event w/in 24 M',
    '# Dept Store Trades w/ MajDerog Event w/in 24 M',
    '# Mortgage Trades w/ MajDerog Event w/in 24 M',
    '# Retail Trades w/ MajDerog Event w/in 24 M',
    '# Sales Finance Trades w/ MajDerog Event w/in 24 M',
    '# Dept Store Trades MajDerog',
    '# Auto Trades MajDerog',
    '# Mortgage Trades MajDerog',
    '# Retail Trades MajDerog',
    '# Dept Store Trades w/ Unpaid MajDerog Event w/in 24 M',
    '# Retail Trades w/ Unpaid MajDerog Event w/in 24 M',
    '# Sales Finance Trades w/ Unpaid MajDerog Event w/in 24 M',
    '# Dept Store Trades Unpaid MajDerog ',
    '# Retail Trades Unpaid MajDerog ',
    '# Open Retail Trades w/ Update w/in 3 M w/ Bal >= 50% High Credit',
    '# Open Auto Trades w/ Update w/in 3 M w/ Bal >= 75% Loan Amount ',
    '# Open Bkcrd Trades w/ Update w/in 3 M w/ Bal >= 75% High Credit ',
    '# Open Dept Store Trades w/ Update w/in 3 M w/ Bal >= 75% High Credit ',
    '# Open Inst Trades w/ Update w/in 3 M w/ Bal >= 75% Loan Amount ',
    '# Open Mortgage Trades w/ Update w/in 3 M w/ Bal >= 75% Loan Amount ',
    '# Retail Trades Reported w/in 3 M',
    '# Retail Trades Reported w/in 6 M',
    '# Dept Store Trades Reported w/in 3 M',
    '# Bkcrd Trades Reported w/in 6 M',
    'Age Newest Date Last Activity Bkcrd Trades Paid as Agreed',
    'Age Newest Date Last Activity Trades Other Than Paid as Agreed',
    'Age Newest Date Last Activity Bkcrd Trades Other Than Paid as Agreed',
    'Age Newest Date Last Activity Inst Trades Other Than Paid as Agreed',
    'Age Newest Date Last Activity Rev Trades Other Than Paid as Agreed',
    'T Collection Amount Unpaid 3rd Party Collections w/in 12 M',
    'T Collection Amount Unpaid 3rd Party  Collections w/in 24 M',
    'T Collection Amount Unpaid 3rd Party Collections',
    'Dismissed Bankruptcy Public Record w/in 24 M Flag',
    'Dismissed Bankruptcy Public Record Flag',
    'Discharged Bankruptcy Public Record w/in 24 M Flag',
    'Discharged Bankruptcy Public Record Flag',
    'Non-Dismissed, Non-Discharged Bankruptcy Public Record w/in 24 M Flag',
    'Non-Dismissed, Non-Discharged Bankruptcy Public Record Flag',
    '# Public Record Tax Liens, Suits and Judgments and 3rd Party Collection Items excluding Debt Buyer',
    'T Bal Defaulted StuLoan Trades w/Update w/in 3 M',
    '% Bal to High Credit Open Retail Trades w/Update w/in 3 M',
    'Age Newest Tax Lien Public Record Item',
    'Age Newest Judgment Public Record Item',
    '# of Tax Lien Pub Rec Item',
    '# of Judgment Pub Rec Item',
    '% Open Dept Store Trades to Open Retail Trades',
    '% Open Dept Store Trades to Dept Store Trades',
    'T Bal Bkcrd Trades w/Update w/in 3 M',
    'T Bal Credit Union Trades w/Update w/in 3 M',
    'T Bal Dept Store Trades w/Update w/in 3 M',
    'T Bal PerFin and StuLoan Trades w/Update w/in 3 M',
    'T Bal Retail Trades w/Update w/in 3 M',
    'T Bal Sales Finance Trades w/Update w/in 3 M',
    '% Bal to High Credit Open Bkcrd Trades w/Update w/in 3 M',
    '% Bal to High Credit Open Dept Store Trades w/Update w/in 3 M',
    '% Bal to T Loan Amount Open Auto Trades w/Update w/in 3 M',
    '% Bal to T Loan Amount Open Inst Trades w/Update w/in 3 M',
    '% Bal to T Loan Amount Open Mortgage Trades w/Update w/in 3 M',
    '% T PD to T Bal Bkcrd Trades w/Update w/in 3 M',
    '% T PD to T Bal Retail Trades w/Update w/in 3 M',
    '% Bkcrd Trades Satis w/in 3 M to Bkcrd Trades Reported w/in 3 M',
    '% Dept Store Trades Satis w/in 3 M to Dept Store Trades Reported w/in 3 M',
    '% Inst Trades Satis w/in 3 M to Inst Trades Reported w/in 3 M',
    '# Dept Store Tr Wo Rat 60 DPD or Worse w/in 3 M or MajDerog Event w/in 24 M',
    '# Dept Store Tr Wo Rat 90 DPD or Worse w/in 3 M or MajDerog Event w/in 24 M',
    '# Dept Store Tr Wo Rat 120-180 or More DPD or Worse w/in 3 M or MajDerog Event w/in 24 M',
    '# Retail Tr Wo Rat 60 DPD or Worse w/in 3 M or MajDerog Event w/in 24 M',
    '# Retail Tr Wo Rat 90 DPD or Worse w/in 3 M or MajDerog Event w/in 24 M',
    '# Retail Tr Wo Rat 120-180 or More DPD or Worse w/in 3 M or MajDerog Event w/in 24 M',
    '# Dept Store Tr Wo Rat No Worse Than 59 DPD w/in 3 M',
    '# Retail Tr Wo Rat No Worse Than 59 DPD w/in 3 M',
    '% Bkcrd Tr Wo Rat 60 DPD or Worse w/in 3 M or MajDerog Event w/in 24 M to Bkcrd Trades Reported w/in 3 M',
    '% Bkcrd Tr Wo Rat 90 DPD or Worse w/in 3 M or MajDerog Event w/in 24 M to Bkcrd Trades Reported w/in 3 M',
    '% Bkcrd Tr Wo Rat 120-180 or More DPD or Worse w/in 3 M or MajDerog Event w/in 24 M to Bkcrd Trades Reported w/in 3 M',
    'T Collection Amount 3rd Party Collections w/in 12 M',
    'T Collection Amount 3rd Party Collections w/in 24 M',
    'T Collection Amount 3rd Party Collections',
    '# Dept Store Tr Wo Rat 60 DPD or Worse w/in 6 M or MajDerog Event w/in 24 M',
    '# Dept Store Tr Wo Rat 90 DPD or Worse w/in 6 M or MajDerog Event w/in 24 M',
    '# Dept Store Tr Wo Rat 120-180 or More DPD or Worse w/in 6 M or MajDerog Event w/in 24 M',
    '# Retail Tr Wo Rat 60 DPD or Worse w/in 6 M or MajDerog Event w/in 24 M',
    '# Retail Tr Wo Rat 90 DPD or Worse w/in 6 M or MajDerog Event w/in 24 M',
    '# Retail Tr Wo Rat 120-180 or More DPD or Worse w/in 6 M or MajDerog Event w/in 24 M',
    '% Bkcrd Trades Satis w/in 6 M to Bkcrd Trades Reported w/in 6 M',
    '% Bkcrd Tr Wo Rat 60 DPD or Worse w/in 6 M or MajDerog Event w/in 24 M to Bkcrd Trades Reported w/in 6 M',
    '% Bkcrd Tr Wo Rat 90 DPD or Worse w/in 6 M or MajDerog Event w/in 24 M to Bkcrd Trades Reported w/in 6 M',
    '% Bkcrd Tr Wo Rat 120-180 or More DPD or Worse w/in 6 M or MajDerog Event w/in 24 M to Bkcrd Trades Reported w/in 6 M',
    '# Dept Store Tr Wo Rat Ever 60 DPD or Worse ',
    '# Dept Store Tr Wo Rat Ever 90 DPD or Worse ',
    '# Dept Store Tr Wo Rat Ever 120-180 or More DPD or Worse ',
    '% Dept Store Trades w/MajDerog Event w/in 24 M to Dept Store Trades',
    '% Dept Store Trades w/Unpaid MajDerog Event w/in 24 M to Dept Store Trades',
    '% Inq w/in 3 M to Inq w/in 12 M',
    '% Inq w/in 3 M to Inq w/in 24 M',
    '% Inq w/in 12 M to Inq w/in 24 M'
)

for ($i = 0; $i -lt $col1.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $col1[$i]
    $ws.Cells.Item($row, 2).Value = $col2[$i]
}
